# refactor: word-generator | fix: word-generators & excel-parser
#
# Replicates, via Excel COM automation, the authoring edit that:
#  - makes "variable" the active sheet/tab (was "constant")
#  - fills column A (contract-type marker "A") down rows 3-13 on sheet "variable"
#    (row 2 already carries it) by painting A2's format onto A3:A13
#  - turns the previously-blank row 14 into a new data row (contract type "AT",
#    act number "АТ-2022-02") that otherwise mirrors row 13's C:I values
#  - updates the saved selection / scroll position on both sheets

$wb = $excel.ActiveWorkbook

$wsVar   = $wb.Worksheets.Item("variable")
$wsConst = $wb.Worksheets.Item("constant")

# ---------------------------------------------------------------------------
# 1) Column A ("A" contract-type marker) for rows 3 through 13 -- format is
#    painted down from A2 (which already holds the value/format), then the
#    value is (re)asserted on every cell.
# ---------------------------------------------------------------------------
$wsVar.Range("A2").Copy()
$wsVar.Range("A3:A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($r = 3; $r -le 13; $r++) {
    $wsVar.Cells.Item($r, 1).Value = "A"
}

# ---------------------------------------------------------------------------
# 2) Row 14 -- was an empty spacer row, becomes a full data row cloned from
#    row 13's contract (C:I), but tagged with the new "AT" contract type.
# ---------------------------------------------------------------------------
$wsVar.Rows.Item(14).RowHeight = 83.4

$wsVar.Range("A2").Copy()
$wsVar.Range("A14").PasteSpecial(-4122)
$wsVar.Range("A2").Copy()
$wsVar.Range("B14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsVar.Range("C13:I13").Copy()
$wsVar.Range("C14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsVar.Cells.Item(14, 2).Value = "АТ-2022-02"
$wsVar.Cells.Item(14, 1).Value = "AT"
$wsVar.Cells.Item(14, 3).Value = $wsVar.Cells.Item(13, 3).Value
$wsVar.Cells.Item(14, 4).Value = $wsVar.Cells.Item(13, 4).Value
$wsVar.Cells.Item(14, 5).Value = $wsVar.Cells.Item(13, 5).Value
$wsVar.Cells.Item(14, 6).Value = $wsVar.Cells.Item(13, 6).Value
$wsVar.Cells.Item(14, 7).Value = $wsVar.Cells.Item(13, 7).Value
$wsVar.Cells.Item(14, 8).Value = $wsVar.Cells.Item(13, 8).Value
$wsVar.Cells.Item(14, 9).Value = $wsVar.Cells.Item(13, 9).Value

# ---------------------------------------------------------------------------
# 3) View state: "variable" becomes the active/selected tab, scrolled near
#    the new row, with A14 selected; "constant" loses its own tab-selected /
#    scroll-position flags as a result of no longer being the active sheet.
# ---------------------------------------------------------------------------
$wsVar.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$wsVar.Range("A14").Select()

Write-Output "edit applied"
